# ==========================================================================
# Update gh-pages to output generated at 456a3b4
#
# This script reproduces, via Excel COM automation, the edits made to
# 广州-漫展信息.xlsx:
#   1. "展览"   (sheet1) - "想去人数" (F column) counter bumps for several rows
#   2. "演出"   (sheet2) - a new event row inserted (2024-09-28 StarRocket),
#                           pushing existing rows down by one, plus the
#                           matching F-column counter bumps
#   3. "本地生活" (sheet3) - F column counter bumps
#   4. "全部类型" (sheet4) - same new event row inserted at its chronological
#                           slot (row 26), plus the matching F-column counter
#                           bumps that mirror sheets 1-3
# ==========================================================================

$wb = $excel.ActiveWorkbook

# --------------------------------------------------------------------------
# Helper: bump a set of F-column (numeric) cells on a given worksheet.
# --------------------------------------------------------------------------
function Set-FValues {
    param($ws, $updates)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}

# --------------------------------------------------------------------------
# 1) 展览 (Exhibition) - simple counter bumps, no structural change
# --------------------------------------------------------------------------
$wsExpo = $wb.Worksheets.Item("展览")

$expoUpdates = @{
    2  = 59
    3  = 1024
    5  = 25
    6  = 461
    7  = 745
    8  = 257
    10 = 44
    11 = 416
    13 = 85
    14 = 864
    15 = 118
    16 = 2009
    17 = 496
    18 = 7726
    19 = 567
    20 = 522
    21 = 64
    24 = 231
}
Set-FValues -ws $wsExpo -updates $expoUpdates

# --------------------------------------------------------------------------
# 2) 演出 (Performance) - insert the new StarRocket show at row 9
# --------------------------------------------------------------------------
$wsShow = $wb.Worksheets.Item("演出")

$wsShow.Rows(9).Insert()

$wsShow.Range("A9").Value = 8
$wsShow.Range("A8").Copy()
$wsShow.Range("A9").PasteSpecial(-4122)

$wsShow.Range("B9").Value = "2024-09-28"
$wsShow.Range("C9").Value = "广州·次元闪耀！！偶像与ACG音乐现场 StarRocket呈现"
$wsShow.Range("D9").Value = "革新路124号太古仓码头4号仓 MAO Livehouse 广州(太古仓店)"
$wsShow.Range("E9").Value = "2024.09.28 15:30-09.28 19:30"
$wsShow.Range("F9").Value = 0
$wsShow.Range("G9").Value = 80
$wsShow.Range("H9").Value = "https://show.bilibili.com/platform/detail.html?id=91767"
$wsShow.Range("I9").Value = "//i0.hdslb.com/bfs/openplatform/202408/82yP4MrJ1724851263539.png"

# Re-sequence the A column (it is a plain 0-based row counter, independent
# of content) for every row pushed down by the insert.
for ($r = 10; $r -le 18; $r++) {
    $wsShow.Range("A$r").Value = $r - 1
}

# --------------------------------------------------------------------------
# 3) 本地生活 (Local Life) - simple counter bumps
# --------------------------------------------------------------------------
$wsLocal = $wb.Worksheets.Item("本地生活")

$localUpdates = @{
    2 = 5552
    3 = 404
    4 = 393
}
Set-FValues -ws $wsLocal -updates $localUpdates

# --------------------------------------------------------------------------
# 4) 全部类型 (All Types) - insert the new StarRocket show at row 26
#    (its chronologically-sorted slot across all three source sheets)
# --------------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("全部类型")

$wsAll.Rows(26).Insert()

$wsAll.Range("A26").Value = 25
$wsAll.Range("A25").Copy()
$wsAll.Range("A26").PasteSpecial(-4122)

$wsAll.Range("B26").Value = "2024-09-28"
$wsAll.Range("C26").Value = "广州·次元闪耀！！偶像与ACG音乐现场 StarRocket呈现"
$wsAll.Range("D26").Value = "革新路124号太古仓码头4号仓 MAO Livehouse 广州(太古仓店)"
$wsAll.Range("E26").Value = "2024.09.28 15:30-09.28 19:30"
$wsAll.Range("F26").Value = 0
$wsAll.Range("G26").Value = 80
$wsAll.Range("H26").Value = "https://show.bilibili.com/platform/detail.html?id=91767"
$wsAll.Range("I26").Value = "//i0.hdslb.com/bfs/openplatform/202408/82yP4MrJ1724851263539.png"

for ($r = 27; $r -le 45; $r++) {
    $wsAll.Range("A$r").Value = $r - 1
}

$allUpdates = @{
    2  = 59
    3  = 5552
    4  = 404
    5  = 393
    7  = 1024
    11 = 25
    12 = 461
    13 = 745
    14 = 257
    17 = 44
    18 = 416
    21 = 85
    23 = 864
    24 = 118
    27 = 2009
    28 = 496
    29 = 7726
    32 = 567
    33 = 522
    34 = 64
    38 = 231
}
Set-FValues -ws $wsAll -updates $allUpdates
